$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 157
$ws1.Range("F6").Value = 2715
$ws1.Range("F8").Value = 1598
$ws1.Range("F9").Value = 7349
$ws1.Range("F11").Value = 7511
$ws1.Range("F14").Value = 5952
$ws1.Range("F15").Value = 3210
$ws1.Range("F16").Value = 3575
$ws1.Range("F17").Value = 3
$ws1.Range("F20").Value = 278
$ws1.Range("F23").Value = 2060
$ws1.Range("F24").Value = 111
$ws1.Range("F28").Value = 940
$ws1.Range("F30").Value = 2562
$ws1.Range("F32").Value = 3120
$ws1.Range("F35").Value = 219
$ws1.Range("F37").Value = 458
$ws1.Range("F38").Value = 1200

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 50

# Sheet 4: 全部类型 (All Types - combined view)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 157
$ws4.Range("F7").Value = 50
$ws4.Range("F9").Value = 2715
$ws4.Range("F10").Value = 1598
$ws4.Range("F15").Value = 7349
$ws4.Range("F17").Value = 7511
$ws4.Range("F20").Value = 5952
$ws4.Range("F21").Value = 3211
$ws4.Range("F22").Value = 3575
$ws4.Range("F24").Value = 278
$ws4.Range("F29").Value = 2060
$ws4.Range("F35").Value = 940
$ws4.Range("F37").Value = 2562
$ws4.Range("F40").Value = 3120
$ws4.Range("F42").Value = 219
$ws4.Range("F45").Value = 458
$ws4.Range("F46").Value = 1200
